# Add a new "Case locations and outbreaks / public exposure site" entry
# (Brighton restaurant) as the newest row, and flag the previously-listed
# sites as no longer the newest ("old" instead of "new").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row right under the header for the new Brighton entry,
# pushing Hallam/Moorabbin/Mordialloc/Wonthaggi down by one row each.
$ws.Rows.Item(2).Insert()
# The inserted row inherits the bold/bordered header formatting from row 1;
# reset it back to plain data-row formatting.
$ws.Range("A2:E2").ClearFormats()

$ws.Range("A2").Value = "Brighton"
$ws.Range("B2").Value = "Sons of Mary Restaurant  14 Spink St, Brighton VIC 3186"
$ws.Range("C2").Value = "24/12/20 10:00am-11:05am"
$ws.Range("D2").Value = "Case ate at restaurant"
$ws.Range("E2").Value = "old"

# The other (now older) entries are no longer the newest addition.
$ws.Range("E3").Value = "old"
$ws.Range("E4").Value = "old"
$ws.Range("E5").Value = "old"
$ws.Range("E6").Value = "old"

# Select the full A:E columns, matching the saved view state.
$ws.Range("A1:E1048576").Select()
